$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5: Drawing Visual - Lines/FPS go from 6 to 18
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 18

# Row 6: remove "Immediate" label from B6 (DrawingVisual now performs like DrawingContext)
$ws.Range("B6").Clear()

# Row 6: Drawing Canvas - Lines/FPS go from 13 to 18
$ws.Range("F6").Value = 18
$ws.Range("G6").Value = 18

# Row 7: Stream Geometry - Lines/FPS go from 0.1 to 0.5
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.5

# Update selected cell to B6 as in the saved file
$ws.Range("B6").Select()
